$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2296511627906977
$ws.Range("C2").Value = 0.4767441860465116
$ws.Range("J2").Value = 0.01453488372093023
$ws.Range("P2").Value = 0.1598837209302326
$ws.Range("S2").Value = 0.1191860465116279
$ws.Range("B3").Value = 0.00558659217877095
$ws.Range("C3").Value = 0.0670391061452514
$ws.Range("J3").Value = 0.0223463687150838
$ws.Range("P3").Value = 0.7206703910614525
$ws.Range("S3").Value = 0.1843575418994413
$ws.Range("J4").Value = 0.03389830508474576
$ws.Range("P4").Value = 0.559322033898305
$ws.Range("S4").Value = 0.4067796610169492
$ws.Range("B6").Value = 0.03149606299212598
$ws.Range("D6").Value = 0.01574803149606299
$ws.Range("F6").Value = 0.06692913385826772
$ws.Range("J6").Value = 0.3110236220472441
$ws.Range("O6").Value = 0.02362204724409449
$ws.Range("Q6").Value = 0.1653543307086614
$ws.Range("R6").Value = 0.05118110236220472
$ws.Range("S6").Value = 0.3346456692913386
$ws.Range("B7").Value = 0.1068376068376068
$ws.Range("D7").Value = 0.01282051282051282
$ws.Range("F7").Value = 0.03846153846153846
$ws.Range("J7").Value = 0.1495726495726496
$ws.Range("O7").Value = 0.02991452991452992
$ws.Range("Q7").Value = 0.1752136752136752
$ws.Range("R7").Value = 0.07692307692307693
$ws.Range("S7").Value = 0.4102564102564102
$ws.Range("B8").Value = 0.09210526315789473
$ws.Range("D8").Value = 0.01691729323308271
$ws.Range("E8").Value = 0.001879699248120301
$ws.Range("F8").Value = 0.05451127819548872
$ws.Range("J8").Value = 0.06954887218045112
$ws.Range("O8").Value = 0.01503759398496241
$ws.Range("Q8").Value = 0.1823308270676692
$ws.Range("R8").Value = 0.1052631578947368
$ws.Range("S8").Value = 0.462406015037594
$ws.Range("B9").Value = 0.1081081081081081
$ws.Range("D9").Value = 0.01801801801801802
$ws.Range("F9").Value = 0.09009009009009009
$ws.Range("J9").Value = 0.05855855855855856
$ws.Range("O9").Value = 0.004504504504504504
$ws.Range("Q9").Value = 0.2207207207207207
$ws.Range("R9").Value = 0.06306306306306306
$ws.Range("S9").Value = 0.4369369369369369
$ws.Range("B10").Value = 0.1015625
$ws.Range("D10").Value = 0.02604166666666667
$ws.Range("E10").Value = 0.0006510416666666666
$ws.Range("F10").Value = 0.05338541666666666
$ws.Range("J10").Value = 0.083984375
$ws.Range("O10").Value = 0.01236979166666667
$ws.Range("Q10").Value = 0.2174479166666667
$ws.Range("R10").Value = 0.09895833333333333
$ws.Range("S10").Value = 0.4055989583333333
$ws.Range("G11").Value = 0.1725
$ws.Range("J11").Value = 0.0825
$ws.Range("K11").Value = 0.24
$ws.Range("L11").Value = 0.49
$ws.Range("S11").Value = 0.015
$ws.Range("G12").Value = 0.6926605504587156
$ws.Range("J12").Value = 0.1467889908256881
$ws.Range("K12").Value = 0.02293577981651376
$ws.Range("L12").Value = 0.09174311926605505
$ws.Range("S12").Value = 0.04587155963302753
$ws.Range("G13").Value = 0.5319148936170213
$ws.Range("J13").Value = 0.425531914893617
$ws.Range("S13").Value = 0.0425531914893617
$ws.Range("F15").Value = 0.04072398190045249
$ws.Range("H15").Value = 0.1357466063348416
$ws.Range("I15").Value = 0.07239819004524888
$ws.Range("J15").Value = 0.3891402714932127
$ws.Range("K15").Value = 0.08144796380090498
$ws.Range("M15").Value = 0.01357466063348416
$ws.Range("O15").Value = 0.03167420814479638
$ws.Range("S15").Value = 0.2352941176470588
$ws.Range("F16").Value = 0.02369668246445497
$ws.Range("H16").Value = 0.1469194312796208
$ws.Range("I16").Value = 0.07109004739336493
$ws.Range("J16").Value = 0.4407582938388626
$ws.Range("K16").Value = 0.1042654028436019
$ws.Range("M16").Value = 0.02369668246445497
$ws.Range("O16").Value = 0.05213270142180094
$ws.Range("S16").Value = 0.1374407582938389
$ws.Range("F17").Value = 0.02862254025044723
$ws.Range("H17").Value = 0.1806797853309481
$ws.Range("I17").Value = 0.07871198568872988
$ws.Range("J17").Value = 0.4490161001788909
$ws.Range("K17").Value = 0.06618962432915922
$ws.Range("M17").Value = 0.02504472271914132
$ws.Range("N17").Value = 0.001788908765652952
$ws.Range("O17").Value = 0.04114490161001789
$ws.Range("S17").Value = 0.1288014311270125
$ws.Range("F18").Value = 0.0199203187250996
$ws.Range("H18").Value = 0.199203187250996
$ws.Range("I18").Value = 0.07171314741035857
$ws.Range("J18").Value = 0.4462151394422311
$ws.Range("K18").Value = 0.08366533864541832
$ws.Range("M18").Value = 0.0199203187250996
$ws.Range("O18").Value = 0.05179282868525897
$ws.Range("S18").Value = 0.1075697211155379
$ws.Range("F19").Value = 0.02168674698795181
$ws.Range("H19").Value = 0.1921686746987952
$ws.Range("I19").Value = 0.07951807228915662
$ws.Range("J19").Value = 0.3728915662650603
$ws.Range("K19").Value = 0.1192771084337349
$ws.Range("M19").Value = 0.01445783132530121
$ws.Range("O19").Value = 0.05481927710843373
$ws.Range("S19").Value = 0.1451807228915663
